# NIT-9011134308.xlsx edit
# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# For this specific workbook the meaningful changes are:
#   1. A new mora period (2509) is added for the existing worker
#      JAVIER ANTONIO BLANCO DE LA ROSA (a new row in the detail table).
#   2. The "VALOR MORA" total (E11) grows to reflect the new period.
#   3. "Cant. Periodos" (F13) increases from 18 to 19.
#   4. The "Novedad de Ingreso" / "Observaciones" header columns (H15/J15)
#      are swapped.
#   5. The "Periodo Mora" column (E) becomes centered for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new detail row right after the current last data row
#    (row 33) so the table keeps growing downward, pushing the
#    signature block (rows 38-39) down to rows 39-40.
# ------------------------------------------------------------------
$ws.Rows.Item(34).Insert()

# Preserve the "closing" (thicker bottom border) look of the old last
# row by moving its format down onto the freshly inserted row...
$ws.Range("B33:J33").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)

# ...and restore row 33 back to the regular body-row look, since it is
# no longer the last row of the table.
$ws.Range("B32:J32").Copy()
$ws.Range("B33:J33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Populate the new row with the new mora period for the same
#    worker as the row above it.
# ------------------------------------------------------------------
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1007974773"
$ws.Range("D34").Value = "JAVIER ANTONIO BLANCO DE LA ROSA"
$ws.Range("E34").Value = "2509"
$ws.Range("F34").Value = 52000
$ws.Range("G34").Value = 1300000

# ------------------------------------------------------------------
# 3. Center the "Periodo Mora" column for every detail row (including
#    the newly added one).
# ------------------------------------------------------------------
$ws.Range("E16:E34").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# 4. Update the totals that reflect the new period.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 930890
$ws.Range("F13").Value = 19

# ------------------------------------------------------------------
# 5. Swap the "Novedad de Ingreso" / "Observaciones" column headers.
# ------------------------------------------------------------------
$ws.Range("H15").Value = "Observaciones"
$ws.Range("J15").Value = "Novedad de Ingreso"

Write-Host "edit complete"
